$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 324, shifting the
# existing rows 324-405 down to 326-407 (the previous last two rows end up
# at 406-407, matching the new dimension A1:R407).
$ws.Rows("324:325").Insert()

# --- New row 324: Camote / 1a nueva(o) / Peru ---
$ws.Range("A324").Value = 10
$ws.Range("B324").Value = "Vega Modelo de Temuco"
$ws.Range("C324").Value = "La Araucanía"
$ws.Range("D324").Value = 44508
$ws.Range("E324").Value = 9
$ws.Range("F324").Value = 100112045
$ws.Range("G324").Value = "Zapallo"
$ws.Range("H324").Value = "Camote"
$ws.Range("I324").Value = "1a nueva(o)"
$ws.Range("J324").Value = 400
$ws.Range("K324").Value = 800
$ws.Range("L324").Value = 800
$ws.Range("M324").Value = 800
$ws.Range("N324").Value = "$/kilo (volumen en unidades)"
$ws.Range("O324").Value = "Perú"
$ws.Range("P324").Value = 800
$ws.Range("Q324").Value = 1
$ws.Range("R324").Value = "Hortaliza"

# --- New row 325: Paine / 1a (guarda) / Región del Maule ---
$ws.Range("A325").Value = 10
$ws.Range("B325").Value = "Vega Modelo de Temuco"
$ws.Range("C325").Value = "La Araucanía"
$ws.Range("D325").Value = 44508
$ws.Range("E325").Value = 9
$ws.Range("F325").Value = 100112045
$ws.Range("G325").Value = "Zapallo"
$ws.Range("H325").Value = "Paine"
$ws.Range("I325").Value = "1a (guarda)"
$ws.Range("J325").Value = 2200
$ws.Range("K325").Value = 300
$ws.Range("L325").Value = 350
$ws.Range("M325").Value = 323
$ws.Range("N325").Value = "$/kilo (volumen en unidades)"
$ws.Range("O325").Value = "Región del Maule"
$ws.Range("P325").Value = 323
$ws.Range("Q325").Value = 1
$ws.Range("R325").Value = "Hortaliza"
